$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12 - "Attempting to generate second Maze"
$ws.Range("A12").Value = "Attempting to generate second Maze"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 45079
$ws.Range("D12").Value = "Tried adding kruskal's algoritm to generate the maze, ended up running into too many issues to fix within scope"
$ws.Range("E12").Value = "X"

# Row 13 - "Attempting to made the player able to go trough the wall if they're the same color"
$ws.Range("A13").Value = "Attempting to made the player able to go trough the wall if they're the same color"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 45079
$ws.Range("D13").Value = "Attempted to fix the Color Maze game to have a unique mechanic of having the player go trough the wall but ended up having too many issues with movement. The objects can assign colors"
$ws.Range("E13").Value = "X"

# Row 14 - "Added limiters to the fields"
$ws.Range("A14").Value = "Added limiters to the fields"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 45079
$ws.Range("D14").Value = "Made it so there is a minimum and a maximum to maze size. Also made them changeable in the unity editor"

# Update the active view/selection to match the new log entries
$ws.Activate()
$ws.Range("D15").Select()
